$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Restyle the data block (9 worker/period rows instead of 7):
#    - old last row (22) carries the "bottom border" styles (21-26)
#      and must become a "middle" row (copy format from row 21).
#    - two brand-new rows (23, 24) are required; row 24 becomes the
#      new bottom row (copy format from the CURRENT row 22, before it
#      gets overwritten), row 23 is a middle row (copy from row 21).
# ------------------------------------------------------------------

# New bottom row 24 <- current row22 formatting (has the bottom border)
$ws.Range("B22:J22").Copy()
$ws.Range("B24:J24").PasteSpecial(-4122)  # xlPasteFormats

# Turn old row22 into a normal middle row + populate new row23 the same way
$ws.Range("B21:J21").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B21:J21").Copy()
$ws.Range("B23:J23").PasteSpecial(-4122)  # xlPasteFormats

# ------------------------------------------------------------------
# 2) Write the new worker / period / value data (rows 16-24)
# ------------------------------------------------------------------
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1002497743"
$ws.Range("D16").Value = "IVAN DAVID ATENCIA DELGADO"
$ws.Range("E16").Value = "2506"
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 1423500

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1002497743"
$ws.Range("D17").Value = "IVAN DAVID ATENCIA DELGADO"
$ws.Range("E17").Value = "2505"
$ws.Range("F17").Value = 22776
$ws.Range("G17").Value = 1423500

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1116043174"
$ws.Range("D18").Value = "JOSE FERNANDO GAITAN GAITAN"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 180000
$ws.Range("G18").Value = 4500000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1116043174"
$ws.Range("D19").Value = "JOSE FERNANDO GAITAN GAITAN"
$ws.Range("E19").Value = "2506"
$ws.Range("F19").Value = 180000
$ws.Range("G19").Value = 4500000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1116043174"
$ws.Range("D20").Value = "JOSE FERNANDO GAITAN GAITAN"
$ws.Range("E20").Value = "2505"
$ws.Range("F20").Value = 180000
$ws.Range("G20").Value = 4500000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1052984539"
$ws.Range("D21").Value = "CARLOS JAVIER JARABA GUTIERREZ"
$ws.Range("E21").Value = "2503"
$ws.Range("F21").Value = 32266
$ws.Range("G21").Value = 1423500

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1052984539"
$ws.Range("D22").Value = "CARLOS JAVIER JARABA GUTIERREZ"
$ws.Range("E22").Value = "2503"
$ws.Range("F22").Value = 49348
$ws.Range("G22").Value = 1423500

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1052984539"
$ws.Range("D23").Value = "CARLOS JAVIER JARABA GUTIERREZ"
$ws.Range("E23").Value = "2502"
$ws.Range("F23").Value = 26572
$ws.Range("G23").Value = 1423500

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1002493246"
$ws.Range("D24").Value = "KEVIN ALBERTO ECHEVERRIA BELEÃ?O"
$ws.Range("E24").Value = "2503"
$ws.Range("F24").Value = 49348
$ws.Range("G24").Value = 1423500

# ------------------------------------------------------------------
# 3) Summary fields above the table
# ------------------------------------------------------------------
$ws.Range("E11").Value = 777250      # VALOR MORA total
$ws.Range("C13").Value = 4           # Cant. Trabajadores
$ws.Range("F13").Value = 5           # Cant. Periodos

# ------------------------------------------------------------------
# 4) Move the signature block down two rows (27/28 -> 29/30)
# ------------------------------------------------------------------
$ws.Range("B27:C28").Copy()
$ws.Range("B29:C30").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H27:J28").Copy()
$ws.Range("H29:J30").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B29").Value = "___________________________________"
$ws.Range("H29").Value = "___________________________________"
$ws.Range("B30").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H30").Value = "FIRMA DEL REPRESENTANTE LEGAL"

$ws.Range("B27:C28").UnMerge()
$ws.Range("H27:J28").UnMerge()
$ws.Range("B27:C28").Clear()
$ws.Range("H27:J28").Clear()

$ws.Range("B29:C29").Merge()
$ws.Range("B30:C30").Merge()
$ws.Range("H29:J29").Merge()
$ws.Range("H30:J30").Merge()
